# Apply the edits described by the commit "Echt wirklich aller aller letzte Korrektur."
# Target sheet: "Übung1" (first worksheet) which contains the network-plan (Netzplan) exercise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Change the number of workers/resources for tasks A and B (column I)
#    I2: 3 -> 1, I3: 2 -> 1
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1

# 2. Two cells that used to hold hard-coded numbers now become live formulas
#    that mirror the pattern used by their neighboring cells.
$ws.Range("P18").Formula = "=R18-P17"
$ws.Range("L27").Formula = "=M28-M25"

# 3. Update the active selection to Q27 (cosmetic, matches the saved view state)
$ws.Activate()
$ws.Range("Q27").Select()

$wb.Save()
